$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Remove the leading "PROJECT / COMPASS" row. Excel's native row-delete
#    shifts every row below it up by one, and (crucially) keeps the tables,
#    merged cells and the sheet dimension in sync automatically.
# ---------------------------------------------------------------------------
$ws.Rows(1).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. A handful of "bullet" rows used to store a literal "-" in column A and
#    the real text in column B. They are now collapsed onto a single cell in
#    column A (same text, B blanked out). Re-use column B's existing format
#    (font color / no special alignment) for column A via a format-only
#    paste, then write the final text and clear column B.
# ---------------------------------------------------------------------------
function Collapse-BulletRow {
    param([string]$RowRef, [string]$Text)

    $colA = "A" + $RowRef
    $colB = "B" + $RowRef

    $ws.Range($colB).Copy() | Out-Null
    $ws.Range($colA).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false

    $ws.Range($colA).Value2 = $Text
    $ws.Range($colB).ClearContents() | Out-Null
}

Collapse-BulletRow "10" "Development of 2 tickets (mention ticket IDs if possible) "
Collapse-BulletRow "11" "Knowledge transfer and documentation updates "
Collapse-BulletRow "42" "Start to prepare release on PROD environment - August 2024 "
Collapse-BulletRow "45" "Write the new US related to Top Transaction Banking | Select PNL Item to aplly for step 2 VER BY BL "
Collapse-BulletRow "46" "Creation of US for Uplouad of Rules"
Collapse-BulletRow "49" "Brainstorm on Monitoring and Performance topics "

# ---------------------------------------------------------------------------
# 3. Two notes that used to be prefixed with "* " lose that prefix (format
#    of the cell itself does not change).
# ---------------------------------------------------------------------------
$ws.Range("A28").Value2 = "For the scope of Sprint 29 we've raised the following priorities: "
$ws.Range("A36").Value2 = "Development of 2 tickets (mention ticket IDs if possible) "

# ---------------------------------------------------------------------------
# 4. Restore the view state recorded in the saved workbook (selected cell).
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("B42").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
